# Generate Report for Handoff
# Adds a new tracked file (176fb172-10f6-4c1e-8de4-4255f7a8c9b2) that has
# reached "Ready for handoff" status, inserting a row for it ahead of the
# existing 256cfca7-bd4a-4052-9e6b-203549a09750 row (alphabetical order) on
# all three worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$HYPERLINK_COLOR = 15570276  # OLE BGR encoding of RGB(0x64,0x95,0xED) == style FF6495ED

function Style-AsHyperlink($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
    $range.Font.Underline = $true
    $range.Font.Color = $HYPERLINK_COLOR
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Hyperlinks don't automatically follow a row insert, so drop them all and
# re-create them once every row is in its final place.
$ws1.Hyperlinks.Delete()

$ws1.Rows.Item(6).Insert()

$ws1.Range("A6").Value = "176fb172-10f6-4c1e-8de4-4255f7a8c9b2.md"
$ws1.Range("B6").Value = "Ready for handoff"
$ws1.Range("C6").Value = "Ready for handoff"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cbacac13fdeb6f9c52ff93e5ab0d984437d4725e/e2e/03f15f77-30da-41b6-8aea-bdc6ce9da5b4.md", "", "", "03f15f77-30da-41b6-8aea-bdc6ce9da5b4.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d62e2ed5d707f1a91367ad185761a3309b73c8bd/e2e/31ff6b21-39a5-440d-8b43-c19aceccf2b2.md", "", "", "31ff6b21-39a5-440d-8b43-c19aceccf2b2.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5fac09d3faabaafcdb0fcb740baef32d72e7c393/e2e/56e27cc1-b2bf-4a3a-a632-2fe9cb1be70a.md", "", "", "56e27cc1-b2bf-4a3a-a632-2fe9cb1be70a.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d62e2ed5d707f1a91367ad185761a3309b73c8bd/e2e/dfc43f13-b8b4-4931-a11c-9654dd1f8b83.md", "", "", "dfc43f13-b8b4-4931-a11c-9654dd1f8b83.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/acf8a8c1c792d0557e542e929322343e750a3470/e2e/176fb172-10f6-4c1e-8de4-4255f7a8c9b2.md", "", "", "176fb172-10f6-4c1e-8de4-4255f7a8c9b2.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/94ed0ee46bacc8643018fd8868014874e057d987/e2e/256cfca7-bd4a-4052-9e6b-203549a09750.md", "", "", "256cfca7-bd4a-4052-9e6b-203549a09750.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/cbacac13fdeb6f9c52ff93e5ab0d984437d4725e/.localization-config", "", "", ".localization-config") | Out-Null

foreach ($addr in @("A2","A3","A4","A5","A6","A7","A8")) {
    Style-AsHyperlink $ws1.Range($addr)
}

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Rows.Item(6).Insert()

$ws2.Range("A6").Value = "176fb172-10f6-4c1e-8de4-4255f7a8c9b2.md"
$ws2.Range("B6").Value = "Ready for handoff"
$ws2.Range("C6").Value = "176fb172-10f6-4c1e-8de4-4255f7a8c9b2.acf8a8c1c792d0557e542e929322343e750a3470.zh-cn.xlf"
$ws2.Range("D6").Value = "2016-03-08 02:02:20"
$ws2.Range("G6").Value = "0001-01-01 00:00:00"
$ws2.Range("H6").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cbacac13fdeb6f9c52ff93e5ab0d984437d4725e/e2e/03f15f77-30da-41b6-8aea-bdc6ce9da5b4.md", "", "", "03f15f77-30da-41b6-8aea-bdc6ce9da5b4.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/025a0ff9cfb239a813cfe3a372fd27f955eaaaa7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/03f15f77-30da-41b6-8aea-bdc6ce9da5b4.8b8ffa16ff065c6359911c854b8e82c908a5754b.zh-cn.xlf", "", "", "03f15f77-30da-41b6-8aea-bdc6ce9da5b4.8b8ffa16ff065c6359911c854b8e82c908a5754b.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ff5afff5c2e0414047ae2b843fb2fbc5a1772742/e2e/03f15f77-30da-41b6-8aea-bdc6ce9da5b4.md", "", "", "03f15f77-30da-41b6-8aea-bdc6ce9da5b4.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6e20883aa23244fced825e1813102bf8139026ec/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/03f15f77-30da-41b6-8aea-bdc6ce9da5b4.8b8ffa16ff065c6359911c854b8e82c908a5754b.zh-cn.xlf", "", "", "03f15f77-30da-41b6-8aea-bdc6ce9da5b4.8b8ffa16ff065c6359911c854b8e82c908a5754b.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d62e2ed5d707f1a91367ad185761a3309b73c8bd/e2e/31ff6b21-39a5-440d-8b43-c19aceccf2b2.md", "", "", "31ff6b21-39a5-440d-8b43-c19aceccf2b2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/568bcbc28d919c68f03a9afeecbe222c1b703c09/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/31ff6b21-39a5-440d-8b43-c19aceccf2b2.2fd36b9d3d702f71696533ff290a1e7d2efb7467.zh-cn.xlf", "", "", "31ff6b21-39a5-440d-8b43-c19aceccf2b2.2fd36b9d3d702f71696533ff290a1e7d2efb7467.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5fac09d3faabaafcdb0fcb740baef32d72e7c393/e2e/56e27cc1-b2bf-4a3a-a632-2fe9cb1be70a.md", "", "", "56e27cc1-b2bf-4a3a-a632-2fe9cb1be70a.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a29776bd0e00eeed7f77ff128bec4465638c30c4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/56e27cc1-b2bf-4a3a-a632-2fe9cb1be70a.763590c2cdb50a8067cd814f07b4245ef5a0b9d4.zh-cn.xlf", "", "", "56e27cc1-b2bf-4a3a-a632-2fe9cb1be70a.763590c2cdb50a8067cd814f07b4245ef5a0b9d4.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d62e2ed5d707f1a91367ad185761a3309b73c8bd/e2e/dfc43f13-b8b4-4931-a11c-9654dd1f8b83.md", "", "", "dfc43f13-b8b4-4931-a11c-9654dd1f8b83.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/568bcbc28d919c68f03a9afeecbe222c1b703c09/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/dfc43f13-b8b4-4931-a11c-9654dd1f8b83.b233ffe110c1a9f06550f22d23ad0bb934b8ddb8.zh-cn.xlf", "", "", "dfc43f13-b8b4-4931-a11c-9654dd1f8b83.b233ffe110c1a9f06550f22d23ad0bb934b8ddb8.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/acf8a8c1c792d0557e542e929322343e750a3470/e2e/176fb172-10f6-4c1e-8de4-4255f7a8c9b2.md", "", "", "176fb172-10f6-4c1e-8de4-4255f7a8c9b2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/acf8a8c1c792d0557e542e929322343e750a3470/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/176fb172-10f6-4c1e-8de4-4255f7a8c9b2.acf8a8c1c792d0557e542e929322343e750a3470.zh-cn.xlf", "", "", "176fb172-10f6-4c1e-8de4-4255f7a8c9b2.acf8a8c1c792d0557e542e929322343e750a3470.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/94ed0ee46bacc8643018fd8868014874e057d987/e2e/256cfca7-bd4a-4052-9e6b-203549a09750.md", "", "", "256cfca7-bd4a-4052-9e6b-203549a09750.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0bb2e751edaad23d80113232974ea6787809fe3f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/256cfca7-bd4a-4052-9e6b-203549a09750.c4aafc9462aa7278eca6a71317def1367cbdbdf7.zh-cn.xlf", "", "", "256cfca7-bd4a-4052-9e6b-203549a09750.c4aafc9462aa7278eca6a71317def1367cbdbdf7.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/cbacac13fdeb6f9c52ff93e5ab0d984437d4725e/.localization-config", "", "", ".localization-config") | Out-Null

foreach ($addr in @("A2","C2","E2","F2","A3","C3","A4","C4","A5","C5","A6","C6","A7","C7","A8")) {
    Style-AsHyperlink $ws2.Range($addr)
}

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Rows.Item(6).Insert()

$ws3.Range("A6").Value = "176fb172-10f6-4c1e-8de4-4255f7a8c9b2.md"
$ws3.Range("B6").Value = "Ready for handoff"
$ws3.Range("C6").Value = "176fb172-10f6-4c1e-8de4-4255f7a8c9b2.acf8a8c1c792d0557e542e929322343e750a3470.de-de.xlf"
$ws3.Range("D6").Value = "2016-03-08 02:02:27"
$ws3.Range("G6").Value = "0001-01-01 00:00:00"
$ws3.Range("H6").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cbacac13fdeb6f9c52ff93e5ab0d984437d4725e/e2e/03f15f77-30da-41b6-8aea-bdc6ce9da5b4.md", "", "", "03f15f77-30da-41b6-8aea-bdc6ce9da5b4.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cfe6beaa57e96521be22d49b9d5463d5fbda2ea4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/03f15f77-30da-41b6-8aea-bdc6ce9da5b4.8b8ffa16ff065c6359911c854b8e82c908a5754b.de-de.xlf", "", "", "03f15f77-30da-41b6-8aea-bdc6ce9da5b4.8b8ffa16ff065c6359911c854b8e82c908a5754b.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8aa9c6d5319fc0f8c2a2ea59a136acc56b79aeae/e2e/03f15f77-30da-41b6-8aea-bdc6ce9da5b4.md", "", "", "03f15f77-30da-41b6-8aea-bdc6ce9da5b4.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f2879195d6c5408962c72a1b06f28010ab2137b6/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/03f15f77-30da-41b6-8aea-bdc6ce9da5b4.8b8ffa16ff065c6359911c854b8e82c908a5754b.de-de.xlf", "", "", "03f15f77-30da-41b6-8aea-bdc6ce9da5b4.8b8ffa16ff065c6359911c854b8e82c908a5754b.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d62e2ed5d707f1a91367ad185761a3309b73c8bd/e2e/31ff6b21-39a5-440d-8b43-c19aceccf2b2.md", "", "", "31ff6b21-39a5-440d-8b43-c19aceccf2b2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/27871dac90f9f3d43a8e9d33b76f7be96c8ae359/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/31ff6b21-39a5-440d-8b43-c19aceccf2b2.2fd36b9d3d702f71696533ff290a1e7d2efb7467.de-de.xlf", "", "", "31ff6b21-39a5-440d-8b43-c19aceccf2b2.2fd36b9d3d702f71696533ff290a1e7d2efb7467.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5fac09d3faabaafcdb0fcb740baef32d72e7c393/e2e/56e27cc1-b2bf-4a3a-a632-2fe9cb1be70a.md", "", "", "56e27cc1-b2bf-4a3a-a632-2fe9cb1be70a.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e94ffc8d66139487e316e2d45339ce8820508017/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/56e27cc1-b2bf-4a3a-a632-2fe9cb1be70a.763590c2cdb50a8067cd814f07b4245ef5a0b9d4.de-de.xlf", "", "", "56e27cc1-b2bf-4a3a-a632-2fe9cb1be70a.763590c2cdb50a8067cd814f07b4245ef5a0b9d4.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d62e2ed5d707f1a91367ad185761a3309b73c8bd/e2e/dfc43f13-b8b4-4931-a11c-9654dd1f8b83.md", "", "", "dfc43f13-b8b4-4931-a11c-9654dd1f8b83.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/27871dac90f9f3d43a8e9d33b76f7be96c8ae359/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/dfc43f13-b8b4-4931-a11c-9654dd1f8b83.b233ffe110c1a9f06550f22d23ad0bb934b8ddb8.de-de.xlf", "", "", "dfc43f13-b8b4-4931-a11c-9654dd1f8b83.b233ffe110c1a9f06550f22d23ad0bb934b8ddb8.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/acf8a8c1c792d0557e542e929322343e750a3470/e2e/176fb172-10f6-4c1e-8de4-4255f7a8c9b2.md", "", "", "176fb172-10f6-4c1e-8de4-4255f7a8c9b2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/acf8a8c1c792d0557e542e929322343e750a3470/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/176fb172-10f6-4c1e-8de4-4255f7a8c9b2.acf8a8c1c792d0557e542e929322343e750a3470.de-de.xlf", "", "", "176fb172-10f6-4c1e-8de4-4255f7a8c9b2.acf8a8c1c792d0557e542e929322343e750a3470.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/94ed0ee46bacc8643018fd8868014874e057d987/e2e/256cfca7-bd4a-4052-9e6b-203549a09750.md", "", "", "256cfca7-bd4a-4052-9e6b-203549a09750.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b5179c3217c8a26a4a4983ba76b3cb2b85d42a6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/256cfca7-bd4a-4052-9e6b-203549a09750.c4aafc9462aa7278eca6a71317def1367cbdbdf7.de-de.xlf", "", "", "256cfca7-bd4a-4052-9e6b-203549a09750.c4aafc9462aa7278eca6a71317def1367cbdbdf7.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/cbacac13fdeb6f9c52ff93e5ab0d984437d4725e/.localization-config", "", "", ".localization-config") | Out-Null

foreach ($addr in @("A2","C2","E2","F2","A3","C3","A4","C4","A5","C5","A6","C6","A7","C7","A8")) {
    Style-AsHyperlink $ws3.Range($addr)
}

Write-Host "Inserted row for 176fb172-10f6-4c1e-8de4-4255f7a8c9b2 across Overview/zh-cn/de-de sheets."
